# Update slides and dates
# 1) Slide 1 ("Welcome" title slide): "Michaelmas 2023" -> "Hilary 2024"
#    (the "Michaelmas" run is removed, the remaining run's text becomes "Hilary 2024")
# 2) Slide 3 (requirements slide): GitHub URL text updated from the
#    Michaelmas-Term-2023 repo link to the Hilary-Term-2024 repo link.

$p = $ppt.ActivePresentation

function Replace-InRange($textRange, [string]$search, [string]$replacement) {
    # Finds literal $search inside $textRange.Text and replaces just that
    # sub-range (via Characters(start,len)) so that surrounding runs and
    # their formatting are left completely untouched.
    $full = $textRange.Text
    $idx = $full.IndexOf($search)
    if ($idx -lt 0) {
        return $false
    }
    $sub = $textRange.Characters($idx + 1, $search.Length)
    $sub.Text = $replacement
    return $true
}

# --- Slide 1: "Michaelmas 2023" -> "Hilary 2024" ---------------------------
# The original text is made of two runs: "Michaelmas" (run A) and " 2023"
# (run B). The edit drops run A entirely and turns run B's text into
# "Hilary 2024", so run B's (not run A's) formatting must be the one that
# survives. We therefore first delete the "Michaelmas" substring (which
# removes run A), then replace the now-adjacent " 2023" text (still run B)
# with "Hilary 2024".
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -like "*Michaelmas 2023*") {
            Replace-InRange $tr "Michaelmas" ""
            $tr2 = $shape.TextFrame.TextRange
            Replace-InRange $tr2 " 2023" "Hilary 2024"
        }
    }
}

# --- Slide 3: GitHub repo link text -----------------------------------------
$slide3 = $p.Slides.Item(3)
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $shape = $slide3.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -like "*https://github.com/sraorao/MSD_R_course*") {
            # The URL and the "_MT2023" suffix are two separate runs, so fix
            # them up with two separate, targeted substring replacements.
            Replace-InRange $tr "https://github.com/sraorao/MSD_R_course" "https://github.com/sraorao"
            $tr2 = $shape.TextFrame.TextRange
            Replace-InRange $tr2 "_MT2023" "/MSD_R_course_HT2024"
        }
    }
}
